$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Update Sheet1 column H: header + values ---
$ws1.Range("H1").Value = "Conexion"
$ws1.Range("H2").Value = 1
$ws1.Range("H3").Value = 1
$ws1.Range("H4").Value = 1

# --- Add new worksheet "Sensores" after Sheet1 ---
$wsSensores = $wb.Worksheets.Add($null, $ws1)
$wsSensores.Name = "Sensores"

# Header row
$wsSensores.Range("A1").Value = "nombre_usuario"
$wsSensores.Range("B1").Value = "Thinkpower"
$wsSensores.Range("C1").Value = "AOVX GL100"

# Copy the header style (bold + border + centered) from Sheet1's A1 header cell
$ws1.Range("A1").Copy()
$wsSensores.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows
$wsSensores.Range("A2").Value = "Mario"
$wsSensores.Range("B2").Value = 1
$wsSensores.Range("C2").Value = 1

$wsSensores.Range("A3").Value = "luigi"
$wsSensores.Range("B3").Value = 1
$wsSensores.Range("C3").Value = 1

$wsSensores.Range("A4").Value = "peach"
$wsSensores.Range("B4").Value = 1
$wsSensores.Range("C4").Value = 1

# --- Selection state: Sensores shows C5 selected, Sheet1 is the active tab
#     with H5 selected (matches the post-edit cursor position) ---
$wsSensores.Range("C5").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("H5").Select() | Out-Null
